$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("B1").Value = "Definition"
$ws.Range("C1").Value = "Preferred Term"
$ws.Range("D1").Value = "Counts for OLMESARTAN MEDOXOMIL"
$ws.Range("E1").Value = "PRR"
$ws.Range("F1").Value = "ROR"

# Data row (row 2)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = '<a href="http://www.merriam-webster.com/medlineplus/SPRUE-LIKE%20ENTEROPATHY" target="_blank">Definition</a>'
$ws.Range("C2").Value = "SPRUE-LIKE ENTEROPATHY"
$ws.Range("D2").Value = 3282.0
$ws.Range("E2").Value = 1463.77
$ws.Range("F2").Value = 348.69
